$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 41.666668
$ws.Range("I11").Value = 41.666668
$ws.Range("K11").Value = 41.666668
$ws.Range("M11").Value = 98.333332
$ws.Range("H99").Value = 491
$ws.Range("J99").Value = 1000
$ws.Range("L99").Value = 3000
$ws.Range("N99").Value = -5996
$ws.Range("H116").Value = 9906.5
$ws.Range("I116").Value = 9891.5
$ws.Range("J116").Value = 9911.5
$ws.Range("K116").Value = 9891.5
$ws.Range("L116").Value = 9911.5
$ws.Range("M116").Value = -6449.5
$ws.Range("N116").Value = -16795.5
$ws.Range("H135").Value = 1884.0869
$ws.Range("I135").Value = 1924.2727
$ws.Range("J135").Value = 1000
$ws.Range("K135").Value = 17318.4543
$ws.Range("L135").Value = 9000
$ws.Range("M135").Value = -14783.4543
$ws.Range("N135").Value = -14070
$ws.Range("H138").Value = 2095.3064
$ws.Range("I138").Value = 831.5454999999999
$ws.Range("J138").Value = 2367.8823
$ws.Range("K138").Value = 2494.6365
$ws.Range("L138").Value = 7103.646900000001
$ws.Range("M138").Value = 2645.3635
$ws.Range("N138").Value = -17383.6469

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1276.38
$ws.Range("I32").Value = 1038.1666
$ws.Range("K32").Value = 1038.1666
$ws.Range("M32").Value = -751.1666
$ws.Range("H39").Value = 19690.445
$ws.Range("J39").Value = 24020
$ws.Range("L39").Value = 24020
$ws.Range("N39").Value = -25060
$ws.Range("H45").Value = 2069
$ws.Range("I45").Value = 2093.25
$ws.Range("K45").Value = 2093.25
$ws.Range("M45").Value = -1716.25
$ws.Range("H61").Value = 5681.593
$ws.Range("I61").Value = 5118.65
$ws.Range("K61").Value = 5118.65
$ws.Range("M61").Value = -4906.65
$ws.Range("H88").Value = 1882.9445
$ws.Range("I88").Value = 1512.4286
$ws.Range("J88").Value = 2118.7273
$ws.Range("K88").Value = 1512.4286
$ws.Range("L88").Value = 2118.7273
$ws.Range("M88").Value = -1106.4286
$ws.Range("N88").Value = -2930.7273
$ws.Range("H91").Value = 1882.9445
$ws.Range("I91").Value = 1512.4286
$ws.Range("J91").Value = 2118.7273
$ws.Range("K91").Value = 1512.4286
$ws.Range("L91").Value = 2118.7273
$ws.Range("M91").Value = -108.4286
$ws.Range("N91").Value = -4926.7273
$ws.Range("H122").Value = 2353.9143
$ws.Range("I122").Value = 1997.5312
$ws.Range("J122").Value = 6155.3335
$ws.Range("K122").Value = 5992.5936
$ws.Range("L122").Value = 18466.0005
$ws.Range("M122").Value = -3542.5936
$ws.Range("N122").Value = -23366.0005
$ws.Range("H123").Value = 50428.5
$ws.Range("J123").Value = 50428.5
$ws.Range("L123").Value = 50428.5
$ws.Range("N123").Value = -60228.5
$ws.Range("H132").Value = 5893.086
$ws.Range("I132").Value = 6966.273
$ws.Range("K132").Value = 20898.819
$ws.Range("M132").Value = -18368.819
$ws.Range("H136").Value = 5681.593
$ws.Range("I136").Value = 5118.65
$ws.Range("K136").Value = 15355.95
$ws.Range("M136").Value = -12805.95

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("H107").Value = 3650
$ws.Range("I107").Value = 3650
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 3650
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -1730
$ws.Range("M82").Value = ""
$ws.Range("M85").Value = ""
$ws.Range("N107").Value = ""

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 93.181816
$ws.Range("I7").Value = 83.94118
$ws.Range("K7").Value = 83.94118
$ws.Range("M7").Value = 29.05882
$ws.Range("H31").Value = 4450.756
$ws.Range("I31").Value = 1049.4286
$ws.Range("K31").Value = 1049.4286
$ws.Range("M31").Value = -754.4286
$ws.Range("H34").Value = 4450.756
$ws.Range("I34").Value = 1049.4286
$ws.Range("K34").Value = 1049.4286
$ws.Range("M34").Value = -847.4286
$ws.Range("H107").Value = 503.42856
$ws.Range("I107").Value = 494.5
$ws.Range("J107").Value = 507
$ws.Range("K107").Value = 494.5
$ws.Range("L107").Value = 507
$ws.Range("M107").Value = 1425.5
$ws.Range("N107").Value = -4347

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 1249.5
$ws.Range("I51").Value = 500
$ws.Range("J51").Value = 1999
$ws.Range("K51").Value = 1500
$ws.Range("L51").Value = 5997
$ws.Range("M51").Value = -1040
$ws.Range("N51").Value = -6917
$ws.Range("H98").Value = 569.9167
$ws.Range("J98").Value = 657.625
$ws.Range("L98").Value = 1972.875
$ws.Range("N98").Value = -4968.875
$ws.Range("H122").Value = 770
$ws.Range("J122").Value = 817.1429000000001
$ws.Range("L122").Value = 7354.2861
$ws.Range("N122").Value = -12254.2861
$ws.Range("H129").Value = 4003530.8
$ws.Range("I129").Value = 4533.273
$ws.Range("J129").Value = 7145600
$ws.Range("K129").Value = 13599.819
$ws.Range("L129").Value = 21436800
$ws.Range("M129").Value = -8599.819
$ws.Range("N129").Value = -21446800

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 585.9286
$ws.Range("I97").Value = 545.8889
$ws.Range("J97").Value = 658
$ws.Range("K97").Value = 545.8889
$ws.Range("L97").Value = 658
$ws.Range("M97").Value = -49.88890000000004
$ws.Range("N97").Value = -1650
$ws.Range("H126").Value = 15007512
$ws.Range("J126").Value = 19506166
$ws.Range("L126").Value = 58518498
$ws.Range("N126").Value = -58523438
$ws.Range("H132").Value = 8339.854499999999
$ws.Range("I132").Value = 8611.079
$ws.Range("J132").Value = 7733.5884
$ws.Range("K132").Value = 25833.237
$ws.Range("L132").Value = 23200.7652
$ws.Range("M132").Value = -23303.237
$ws.Range("N132").Value = -28260.7652

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5245.3887
$ws.Range("I40").Value = 5318.647
$ws.Range("K40").Value = 5318.647
$ws.Range("M40").Value = -5182.647
$ws.Range("H61").Value = 6761.926
$ws.Range("I61").Value = 6384.048
$ws.Range("J61").Value = 8084.5
$ws.Range("K61").Value = 6384.048
$ws.Range("L61").Value = 8084.5
$ws.Range("M61").Value = -6182.048
$ws.Range("N61").Value = -8488.5
$ws.Range("H68").Value = 3813.8462
$ws.Range("I68").Value = 2816.3333
$ws.Range("J68").Value = 15784
$ws.Range("K68").Value = 2816.3333
$ws.Range("L68").Value = 15784
$ws.Range("M68").Value = -2067.3333
$ws.Range("N68").Value = -17282
$ws.Range("H71").Value = 3813.8462
$ws.Range("I71").Value = 2816.3333
$ws.Range("J71").Value = 15784
$ws.Range("K71").Value = 14081.6665
$ws.Range("L71").Value = 78920
$ws.Range("M71").Value = -10337.6665
$ws.Range("N71").Value = -86408
$ws.Range("H113").Value = 6761.926
$ws.Range("I113").Value = 6384.048
$ws.Range("J113").Value = 8084.5
$ws.Range("K113").Value = 6384.048
$ws.Range("L113").Value = 8084.5
$ws.Range("M113").Value = -4214.048
$ws.Range("N113").Value = -12424.5
$ws.Range("H132").Value = 10000
$ws.Range("J132").Value = 12000
$ws.Range("L132").Value = 36000
$ws.Range("N132").Value = -41060
$ws.Range("H136").Value = 3587.7144
$ws.Range("I136").Value = 2820.1667
$ws.Range("J136").Value = 3894.7334
$ws.Range("K136").Value = 8460.500100000001
$ws.Range("L136").Value = 11684.2002
$ws.Range("M136").Value = -5910.500100000001
$ws.Range("N136").Value = -16784.2002

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3598.9744
$ws.Range("I122").Value = 3496.2163
$ws.Range("K122").Value = 10488.6489
$ws.Range("M122").Value = -8038.6489
$ws.Range("H126").Value = 3939.3333
$ws.Range("I126").Value = 3388.3635
$ws.Range("K126").Value = 10165.0905
$ws.Range("M126").Value = -7695.0905
$ws.Range("H136").Value = 5024.5557
$ws.Range("J136").Value = 7845.143
$ws.Range("L136").Value = 23535.429
$ws.Range("N136").Value = -28635.429
